$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$default = -19.60781069860693

# Fill default value across B2:K21
for ($r = 2; $r -le 21; $r++) {
    for ($c = 2; $c -le 11; $c++) {
        $ws.Cells.Item($r, $c).Value = $default
    }
}

# Overrides (non-default values)
$ws.Cells.Item(2, 3).Value = 1.98163609527907
$ws.Cells.Item(3, 9).Value = 1.366528963626321
$ws.Cells.Item(4, 3).Value = 1.985644711537074
$ws.Cells.Item(4, 4).Value = 1.68664907430642
$ws.Cells.Item(4, 6).Value = 3.432476513723809
$ws.Cells.Item(4, 8).Value = 1.551253137691483
$ws.Cells.Item(5, 3).Value = 1.692806072707717
$ws.Cells.Item(5, 7).Value = 2.857759940812927
$ws.Cells.Item(7, 2).Value = 2.456604802728271
$ws.Cells.Item(8, 5).Value = 1.804912638902015
$ws.Cells.Item(9, 2).Value = 3.859049548533317
$ws.Cells.Item(10, 9).Value = 1.614076693281663
$ws.Cells.Item(10, 11).Value = 2.341494888169733
$ws.Cells.Item(11, 5).Value = 2.916198402861638
$ws.Cells.Item(11, 7).Value = 2.838163212284698
$ws.Cells.Item(11, 11).Value = 1.775858265932455
$ws.Cells.Item(13, 5).Value = 2.513418308760795
$ws.Cells.Item(13, 11).Value = 1.85942318116096
$ws.Cells.Item(14, 4).Value = 1.5471810385086
$ws.Cells.Item(14, 11).Value = 2.027806910307834
$ws.Cells.Item(15, 4).Value = 1.75347612540035
$ws.Cells.Item(17, 3).Value = 2.124155029985084
$ws.Cells.Item(17, 4).Value = 1.810503645510569
$ws.Cells.Item(17, 8).Value = 2.065902521909486
$ws.Cells.Item(17, 9).Value = 2.142183493712011
$ws.Cells.Item(18, 8).Value = 2.066066390483264
$ws.Cells.Item(18, 9).Value = 2.038386297528072
$ws.Cells.Item(18, 10).Value = 4.321926379511744
$ws.Cells.Item(19, 4).Value = 2.053367630034424
$ws.Cells.Item(19, 8).Value = 1.665937671564617
$ws.Cells.Item(19, 9).Value = 1.819167023954747
$ws.Cells.Item(20, 3).Value = 1.036028760598511
$ws.Cells.Item(20, 4).Value = 1.499966372278189
$ws.Cells.Item(20, 6).Value = 3.2021969986196
$ws.Cells.Item(20, 8).Value = 1.566769425201169
$ws.Cells.Item(20, 9).Value = 1.207520212745211
$ws.Cells.Item(20, 11).Value = 1.92631024325964
$ws.Cells.Item(21, 3).Value = 1.291178097743003
$ws.Cells.Item(21, 5).Value = 1.699392890487549
$ws.Cells.Item(21, 7).Value = 2.485426633706928
$ws.Cells.Item(21, 8).Value = 1.355732428449048
